# Auto-generated Excel COM-interop script to apply the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued cells (Coin name, Link URL, Volume label) ---
# These are plain text columns; assigning directly keeps them as text.
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('E6').Value = '5GateTokenGT'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E8').Value = '7MXTokenMX'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E9').Value = '8FTXTokenFTT'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B19').Value = 'UpBots'
$ws.Range('C19').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('E19').Value = '18UpBotsUBXT'
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B24').Value = 'BTSEToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('E24').Value = '23BTSETokenBTSE'
$ws.Range('B25').Value = 'One'
$ws.Range('C25').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E25').Value = '24OneONEBestin24h'
$ws.Range('E27').Value = '26ProBitTokenPROB'

# --- Numeric-looking text cells (Price column D) ---
# These must remain text (not be coerced into real numbers) so that exact
# formatting (e.g. trailing zeros) is preserved, matching the source data.
$priceCells = @(
    'D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50'
)
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '247.08'
$ws.Range('D3').Value = '22.66'
$ws.Range('D4').Value = '5.547'
$ws.Range('D5').Value = '0.05604'
$ws.Range('D6').Value = '3.402'
$ws.Range('D7').Value = '6.468'
$ws.Range('D8').Value = '0.8024'
$ws.Range('D9').Value = '1.056'
$ws.Range('D10').Value = '0.1423'
$ws.Range('D11').Value = '0.07263'
$ws.Range('D12').Value = '0.03193'
$ws.Range('D13').Value = '0.02969'
$ws.Range('D14').Value = '0.09253'
$ws.Range('D15').Value = '0.001660'
$ws.Range('D16').Value = '2.971'
$ws.Range('D17').Value = '0.04691'
$ws.Range('D18').Value = '0.006272'
$ws.Range('D19').Value = '0.007497'
$ws.Range('D20').Value = '0.001050'
$ws.Range('D21').Value = '0.003808'
$ws.Range('D22').Value = '0.0001502'
$ws.Range('D23').Value = '3.981'
$ws.Range('D24').Value = '2.113'
$ws.Range('D25').Value = '0.01162'
$ws.Range('D27').Value = '0.1293'
$ws.Range('D40').Value = '0.04178'
$ws.Range('D41').Value = '0.006960'
$ws.Range('D42').Value = '0.1038'
$ws.Range('D43').Value = '0.003152'
$ws.Range('D44').Value = '0.01018'
$ws.Range('D45').Value = '0.00005633'
$ws.Range('D46').Value = '0.00000000751'
$ws.Range('D47').Value = '0.6810'
$ws.Range('D48').Value = '0.02656'
$ws.Range('D49').Value = '0.00002103'
$ws.Range('D50').Value = '0.01011'

foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "General"
}
